# Rename the "Others" sheet to "Trees" and populate it with the
# NeetCode150 Trees problem list.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Others")
$ws.Name = "Trees"

# Data rows to append below the existing header row (row 1).
# Columns: A=Date Solved, B=Name, C=Algorithm, D=Difficulty,
#          E=Solved First Time, F=Revisit?, G=Understand?
$rows = @(
    @{ A = $null;         B = "Invert Binary Tree";                                         C = "Trees"; D = "Easy";   E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Maximum Depth of Binary Tree";                                C = "Trees"; D = "Easy";   E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Diameter of Binary Tree";                                     C = "Trees"; D = "Easy";   E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Balanced Binary Tree";                                        C = "Trees"; D = "Easy";   E = $null;  F = "Yes"; G = "Yes" },
    @{ A = "08/31/2025";  B = "Same Tree";                                                   C = "Trees"; D = "Easy";   E = "Yes";  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Subtree of Another Tree";                                     C = "Trees"; D = "Easy";   E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Lowest Common Ancestor of a Binary Search Tree";               C = "Trees"; D = "Medium"; E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Binary Tree Level Order Traversal";                           C = "Trees"; D = "Medium"; E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Binary Tree Right Side View";                                 C = "Trees"; D = "Medium"; E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Count Good Nodes in Binary Tree";                             C = "Trees"; D = "Medium"; E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Validate Binary Search Tree";                                 C = "Trees"; D = "Medium"; E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Kth Smallest Element in a Bst";                               C = "Trees"; D = "Medium"; E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Construct Binary Tree From Preorder and Inorder Traversal";    C = "Trees"; D = "Medium"; E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Binary Tree Maximum Path Sum";                                C = "Trees"; D = "Hard";   E = $null;  F = "Yes"; G = "Yes" },
    @{ A = $null;         B = "Serialize and Deserialize Binary Tree";                       C = "Trees"; D = "Hard";   E = $null;  F = "Yes"; G = "Yes" }
)

# Column A holds dates typed as plain text in this workbook (e.g. "08/09/2025"
# in the other sheets), not real Excel date serials, so force text format
# before writing so it round-trips as a string rather than a date number.
$dateCol = $ws.Range("A2:A16")
$dateCol.NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    if ($row.A) { $ws.Cells.Item($r, 1).Value = $row.A }
    if ($row.B) { $ws.Cells.Item($r, 2).Value = $row.B }
    if ($row.C) { $ws.Cells.Item($r, 3).Value = $row.C }
    if ($row.D) { $ws.Cells.Item($r, 4).Value = $row.D }
    if ($row.E) { $ws.Cells.Item($r, 5).Value = $row.E }
    if ($row.F) { $ws.Cells.Item($r, 6).Value = $row.F }
    if ($row.G) { $ws.Cells.Item($r, 7).Value = $row.G }
    $r++
}
